$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data (row 4) appended after the existing rows (1-3)
$ws.Cells.Item(4, 1).Value = 42602.583622685182
$ws.Cells.Item(4, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(4, 2).Value = "Bag"
$ws.Cells.Item(4, 3).Value = 9039
$ws.Cells.Item(4, 4).Value = 10498
$ws.Cells.Item(4, 5).Value = 1277
$ws.Cells.Item(4, 6).Value = 156
$ws.Cells.Item(4, 7).Value = 76
$ws.Cells.Item(4, 8).Value = 67
$ws.Cells.Item(4, 9).Value = 32
$ws.Cells.Item(4, 10).Value = 4
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 57
$ws.Cells.Item(4, 13).Value = 42
